$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-24
$values = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 6, 6),
    @(5, 6, 6),
    @(6, 8, 8),
    @(7, 7, 8),
    @(8, 8, 8),
    @(9, 4, 5),
    @(10, 7, 7),
    @(11, 9, 9),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 5, 5),
    @(21, 4, 4),
    @(22, 6, 6),
    @(23, 4, 4),
    @(24, 3, 3)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
